$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 2768.5
$ws.Range("I92").Value = 1705.619
$ws.Range("K92").Value = 1705.619
$ws.Range("M92").Value = -457.6189999999999
$ws.Range("H96").Value = 1518.3
$ws.Range("I96").Value = 2284.5
$ws.Range("J96").Value = 1007.5
$ws.Range("K96").Value = 6853.5
$ws.Range("L96").Value = 3022.5
$ws.Range("M96").Value = -5480.5
$ws.Range("N96").Value = -5768.5
$ws.Range("H106").Value = 3654.5
$ws.Range("I106").Value = 1881.9166
$ws.Range("K106").Value = 1881.9166
$ws.Range("M106").Value = -1250.9166
$ws.Range("H107").Value = 429
$ws.Range("J107").Value = 400
$ws.Range("L107").Value = 400
$ws.Range("N107").Value = -4240
$ws.Range("H116").Value = 6499.7144
$ws.Range("I116").Value = 6633
$ws.Range("J116").Value = 6399.75
$ws.Range("K116").Value = 6633
$ws.Range("L116").Value = 6399.75
$ws.Range("M116").Value = -3191
$ws.Range("N116").Value = -13283.75
$ws.Range("H132").Value = 1002.2917
$ws.Range("I132").Value = 937.1739
$ws.Range("J132").Value = 2500
$ws.Range("K132").Value = 2811.5217
$ws.Range("L132").Value = 7500
$ws.Range("M132").Value = -281.5217000000002
$ws.Range("N132").Value = -12560
$ws.Range("H137").Value = 3088.484
$ws.Range("I137").Value = 2659.182
$ws.Range("K137").Value = 7977.545999999999
$ws.Range("M137").Value = -5427.545999999999
$ws.Range("H141").Value = 3300.5625
$ws.Range("I141").Value = 3253.9333
$ws.Range("K141").Value = 9761.7999
$ws.Range("M141").Value = -4581.7999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4460.7017
$ws.Range("I32").Value = 3412.2964
$ws.Range("K32").Value = 3412.2964
$ws.Range("M32").Value = -3125.2964
$ws.Range("H45").Value = 3087.7058
$ws.Range("I45").Value = 2490.182
$ws.Range("J45").Value = 4183.1665
$ws.Range("K45").Value = 2490.182
$ws.Range("L45").Value = 4183.1665
$ws.Range("M45").Value = -2113.182
$ws.Range("N45").Value = -4937.1665
$ws.Range("H61").Value = 4056.724
$ws.Range("I61").Value = 3391.3684
$ws.Range("J61").Value = 5320.9
$ws.Range("K61").Value = 3391.3684
$ws.Range("L61").Value = 5320.9
$ws.Range("M61").Value = -3179.3684
$ws.Range("N61").Value = -5744.9
$ws.Range("H74").Value = 19611096
$ws.Range("I74").Value = 22224358
$ws.Range("J74").Value = 11632
$ws.Range("K74").Value = 22224358
$ws.Range("L74").Value = 11632
$ws.Range("M74").Value = -22223484
$ws.Range("N74").Value = -13380
$ws.Range("H77").Value = 19611096
$ws.Range("I77").Value = 22224358
$ws.Range("J77").Value = 11632
$ws.Range("K77").Value = 111121790
$ws.Range("L77").Value = 58160
$ws.Range("M77").Value = -111117422
$ws.Range("N77").Value = -66896
$ws.Range("H97").Value = 1073.3334
$ws.Range("I97").Value = 1073.3334
$ws.Range("K97").Value = 1073.3334
$ws.Range("M97").Value = -577.3334
$ws.Range("H102").Value = 2737
$ws.Range("J102").Value = 7998
$ws.Range("L102").Value = 7998
$ws.Range("N102").Value = -11242
$ws.Range("H132").Value = 2159.8462
$ws.Range("I132").Value = 1622.8206
$ws.Range("J132").Value = 3770.923
$ws.Range("K132").Value = 4868.4618
$ws.Range("L132").Value = 11312.769
$ws.Range("M132").Value = -2338.4618
$ws.Range("N132").Value = -16372.769
$ws.Range("H136").Value = 4056.724
$ws.Range("I136").Value = 3391.3684
$ws.Range("J136").Value = 5320.9
$ws.Range("K136").Value = 10174.1052
$ws.Range("L136").Value = 15962.7
$ws.Range("M136").Value = -7624.1052
$ws.Range("N136").Value = -21062.7

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2327.4546
$ws.Range("I134").Value = 1545.6364
$ws.Range("J134").Value = 3891.0908
$ws.Range("K134").Value = 4636.9092
$ws.Range("L134").Value = 11673.2724
$ws.Range("M134").Value = -2101.9092
$ws.Range("N134").Value = -16743.2724

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 22200.566
$ws.Range("I31").Value = 2059.372
$ws.Range("K31").Value = 2059.372
$ws.Range("M31").Value = -1764.372
$ws.Range("H34").Value = 22200.566
$ws.Range("I34").Value = 2059.372
$ws.Range("K34").Value = 2059.372
$ws.Range("M34").Value = -1857.372
$ws.Range("H58").Value = 3053
$ws.Range("I58").Value = 1438.762
$ws.Range("J58").Value = 7290.375
$ws.Range("K58").Value = 1438.762
$ws.Range("L58").Value = 7290.375
$ws.Range("M58").Value = -1235.762
$ws.Range("N58").Value = -7696.375
$ws.Range("H99").Value = 2206.9167
$ws.Range("I99").Value = 1985.3334
$ws.Range("J99").Value = 2428.5
$ws.Range("K99").Value = 1985.3334
$ws.Range("L99").Value = 2428.5
$ws.Range("M99").Value = -487.3334
$ws.Range("N99").Value = -5424.5
$ws.Range("H126").Value = 2206.9167
$ws.Range("I126").Value = 1985.3334
$ws.Range("J126").Value = 2428.5
$ws.Range("K126").Value = 5956.0002
$ws.Range("L126").Value = 7285.5
$ws.Range("M126").Value = -3486.0002
$ws.Range("N126").Value = -12225.5
$ws.Range("H132").Value = 2726.1555
$ws.Range("I132").Value = 2068.125
$ws.Range("J132").Value = 4345.923
$ws.Range("K132").Value = 6204.375
$ws.Range("L132").Value = 13037.769
$ws.Range("M132").Value = -3674.375
$ws.Range("N132").Value = -18097.769
$ws.Range("H136").Value = 3053
$ws.Range("I136").Value = 1438.762
$ws.Range("J136").Value = 7290.375
$ws.Range("K136").Value = 4316.286
$ws.Range("L136").Value = 21871.125
$ws.Range("M136").Value = -1766.286
$ws.Range("N136").Value = -26971.125

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 11921.75
$ws.Range("J32").Value = 15837.333
$ws.Range("L32").Value = 47511.999
$ws.Range("N32").Value = -48077.999
$ws.Range("H39").Value = 2210.125
$ws.Range("J39").Value = 2530.5
$ws.Range("L39").Value = 7591.5
$ws.Range("N39").Value = -8179.5
$ws.Range("H55").Value = 4099.96
$ws.Range("I55").Value = 1534.4166
$ws.Range("K55").Value = 4603.2498
$ws.Range("M55").Value = -4426.2498
$ws.Range("H114").Value = 4081.5
$ws.Range("I114").Value = 3049.75
$ws.Range("J114").Value = 4494.2
$ws.Range("K114").Value = 9149.25
$ws.Range("L114").Value = 13482.6
$ws.Range("M114").Value = -5895.25
$ws.Range("N114").Value = -19990.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2727.4
$ws.Range("I102").Value = 1515.3077
$ws.Range("K102").Value = 1515.3077
$ws.Range("M102").Value = 106.6922999999999
$ws.Range("H126").Value = 4715.2144
$ws.Range("I126").Value = 2416.5
$ws.Range("J126").Value = 6439.25
$ws.Range("K126").Value = 7249.5
$ws.Range("L126").Value = 19317.75
$ws.Range("M126").Value = -4779.5
$ws.Range("N126").Value = -24257.75
$ws.Range("H132").Value = 2602.5715
$ws.Range("I132").Value = 2213.425
$ws.Range("K132").Value = 6640.275000000001
$ws.Range("M132").Value = -4110.275000000001
$ws.Range("H137").Value = 59971.25
$ws.Range("J137").Value = 59971.25
$ws.Range("L137").Value = 59971.25
$ws.Range("N137").Value = -70171.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5785.864
$ws.Range("I7").Value = 3770.0588
$ws.Range("J7").Value = 12639.6
$ws.Range("K7").Value = 3770.0588
$ws.Range("L7").Value = 12639.6
$ws.Range("M7").Value = -3658.0588
$ws.Range("N7").Value = -12863.6
$ws.Range("H16").Value = 2151.2173
$ws.Range("I16").Value = 499.1579
$ws.Range("K16").Value = 499.1579
$ws.Range("M16").Value = -329.1579
$ws.Range("H22").Value = 2772.4092
$ws.Range("I22").Value = 791.44446
$ws.Range("J22").Value = 4143.846
$ws.Range("K22").Value = 791.44446
$ws.Range("L22").Value = 4143.846
$ws.Range("M22").Value = -496.44446
$ws.Range("N22").Value = -4733.846
$ws.Range("H27").Value = 2772.4092
$ws.Range("I27").Value = 791.44446
$ws.Range("J27").Value = 4143.846
$ws.Range("K27").Value = 791.44446
$ws.Range("L27").Value = 4143.846
$ws.Range("M27").Value = -684.44446
$ws.Range("N27").Value = -4357.846
$ws.Range("H61").Value = 3021.7
$ws.Range("I61").Value = 1997.0667
$ws.Range("J61").Value = 6095.6
$ws.Range("K61").Value = 1997.0667
$ws.Range("L61").Value = 6095.6
$ws.Range("M61").Value = -1795.0667
$ws.Range("N61").Value = -6499.6
$ws.Range("H113").Value = 3021.7
$ws.Range("I113").Value = 1997.0667
$ws.Range("J113").Value = 6095.6
$ws.Range("K113").Value = 1997.0667
$ws.Range("L113").Value = 6095.6
$ws.Range("M113").Value = 172.9332999999999
$ws.Range("N113").Value = -10435.6
$ws.Range("H122").Value = 505999.25
$ws.Range("I122").Value = 670798.5
$ws.Range("K122").Value = 2012395.5
$ws.Range("M122").Value = -2009945.5
$ws.Range("H126").Value = 5785.864
$ws.Range("I126").Value = 3770.0588
$ws.Range("J126").Value = 12639.6
$ws.Range("K126").Value = 11310.1764
$ws.Range("L126").Value = 37918.8
$ws.Range("M126").Value = -8840.1764
$ws.Range("N126").Value = -42858.8
$ws.Range("H136").Value = 3292.4878
$ws.Range("I136").Value = 2437.0938
$ws.Range("J136").Value = 6333.8887
$ws.Range("K136").Value = 7311.2814
$ws.Range("L136").Value = 19001.6661
$ws.Range("M136").Value = -4761.2814
$ws.Range("N136").Value = -24101.6661

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 894.9211
$ws.Range("I100").Value = 712.3
$ws.Range("K100").Value = 1424.6
$ws.Range("M100").Value = -883.5999999999999
$ws.Range("H132").Value = 2346.4443
$ws.Range("I132").Value = 1890.1666
$ws.Range("K132").Value = 5670.4998
$ws.Range("M132").Value = -3140.4998
$ws.Range("H136").Value = 3621.4866
$ws.Range("I136").Value = 2225.2307
$ws.Range("K136").Value = 6675.6921
$ws.Range("M136").Value = -4125.6921
